$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata" ---
$meta = $wb.Worksheets.Item("Metadata")

# Version: 5.0.0 -> 6.0.0
$meta.Range("B3").Value = "6.0.0"

# Date: 2021-12-16T17:36:56+00:00 -> 2022-01-21T20:46:54+00:00
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value: (empty) -> "Alvearie Team"
$meta.Range("B9").Value = "Alvearie Team"

# Row 10 was "Contact" / "No display for ContactDetail" -> becomes "Jurisdiction" / "United States of America"
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# Row 11 was a duplicate "Contact" / "No display for ContactDetail" row that must be removed entirely,
# shifting the remaining rows (old 12-21) up to (new 11-20).
$meta.Rows.Item(11).Delete()

# --- Sheet "Elements" ---
$elements = $wb.Worksheets.Item("Elements")

# Row 2 (the root Extension element): Short/Definition columns (K/L) get the actual
# extension title/description instead of the generic placeholder text.
$elements.Range("K2").Value = "Medical Coverage Indicator"
$elements.Range("L2").Value = "Indicates whether the member has medical benefit coverage: Y or N"
